# ---------------------------------------------------------------------------
# Applies the "Final submission of documentation and fixed the logout button"
# commit to IT.docx:
#   1. Adds a new "kivymd" bullet after the "Google maps API" bullet.
#   2. Extends the functional-testing paragraph with a new sentence about the
#      admin widget, and relocates the "_GoBack" bookmark to sit right before
#      "by the members." (where Word leaves it after the last edit made).
#   3. Removes the old "_GoBack" bookmark from the end of the non-functional
#      testing paragraph (it moved, per step 2).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) New bulleted list item "kivymd" right after the "Google maps API" item.
# ---------------------------------------------------------------------------
$gmaps = $d.Content
$gmaps.Find.Execute("Google maps API", $false) | Out-Null
$gmaps.Collapse(0)
$gmaps.InsertParagraphAfter()

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Google maps API*") {
        $newPara = $p.Next()
        $newPara.Range.Text = "kivymd"
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Relocate the "_GoBack" bookmark: delete it from its old home (end of the
#    "Execution-based Non-Functional Testing" narrative paragraph) ...
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# ... then re-create it inside the functional-testing paragraph, between
#     "... All these functionalities are tested " and "by the members.",
#     and append the new sentence describing the admin widget right after
#     "by the members.".
# ---------------------------------------------------------------------------
$functional = $d.Content
$functional.Find.Execute("tested by the members.", $false) | Out-Null

$bookmarkSpot = $functional.Duplicate
$bookmarkSpot.Collapse(1)
$bookmarkSpot.MoveEnd(1, "tested ".Length)
$bookmarkSpot.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null

$insertSpot = $functional.Duplicate
$insertSpot.Collapse(0)
$insertSpot.InsertAfter(" We created a admin widget which is the admin version of login and card system where we are able to see the events in the database and corresponding information which is then tested to see if it is accurate and runs on the app.")
